$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update H column text values (question "reason" text), keeping the same
# wording except for rows 2 and 5 which are rephrased.
$ws.Range("H2").Value  = "two or more of your family members have been affected by breast cancer"
$ws.Range("H3").Value  = "your father has been affected by breast cancer"
$ws.Range("H4").Value  = "your brother has been affected by breast cancer"
$ws.Range("H5").Value  = "two or more of your family members have been affected by ovarian cancer"
$ws.Range("H6").Value  = "you have family members have affected by both breast and ovarian cancer"
$ws.Range("H8").Value  = "you have more than one family member affected by breast cancer"
$ws.Range("H9").Value  = "your mother was under 60 when she was affected by breast cancer"
$ws.Range("H11").Value = "your sister was under 60 when she was affected by breast cancer"
$ws.Range("H13").Value = "your grandmother was under 40 when she was affected by breast cancer"
$ws.Range("H15").Value = "your half-sister was under 40 when she was affected by breast cancer"
$ws.Range("H17").Value = "your aunt was under 40 when she was affected by breast cancer"
$ws.Range("H19").Value = "your niece was under 40 when she was affected by breast cancer"
$ws.Range("H21").Value = "you have more than one family member affected by ovarian cancer"

# Update C column numeric values (item_identifier counters shift by one for
# most rows, row 7 shifts down by one).
$ws.Range("C7").Value  = 5
$ws.Range("C9").Value  = 2
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = 3
$ws.Range("C12").Value = 3
$ws.Range("C13").Value = 4
$ws.Range("C14").Value = 4
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 5
$ws.Range("C17").Value = 6
$ws.Range("C18").Value = 6
$ws.Range("C19").Value = 7
$ws.Range("C20").Value = 7

# Update the active selection/cell to match the saved view state.
$ws.Range("H11").Select()
